$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3, column A: text changes from admin123@example.com to admin@example.com,
# and the hyperlink display font switches from Roboto 14 (underline) to Arial 14 (underline)
$ws.Range("A3").Value = "admin@example.com"
$ws.Range("A3").Font.Name = "Arial"
$ws.Range("A3").Font.Size = 14

# C2 / C3: font changes from Arial 12 to Arial 14 (keep existing left/center alignment, no fill)
$ws.Range("C2").Font.Size = 14
$ws.Range("C2").HorizontalAlignment = -4131
$ws.Range("C2").VerticalAlignment = -4108

$ws.Range("C3").Font.Size = 14
$ws.Range("C3").HorizontalAlignment = -4131
$ws.Range("C3").VerticalAlignment = -4108

# Update the active selection to match the saved workbook view
$ws.Range("B6").Select()
